$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Moorings")
$ws2 = $wb.Worksheets.Item("Asset_Cal_Info")

# --- Sheet1 (Moorings) corrections ---
# Correct Lat/Lon text formatting (straight apostrophe, no inner spaces)
$ws1.Range("G2").Value = "40°5.801'N"
$ws1.Range("H2").Value = "70°52.764'W"

# Add decimal-degree helper formulas in L2 / M2
$ws1.Range("L2").Formula = '=((LEFT(G2,(FIND("°",G2,1)-1)))+(MID(G2,(FIND("°",G2,1)+1),(FIND("''",G2,1))-(FIND("°",G2,1)+1))/60))*(IF(RIGHT(G2,1)="N",1,-1))'
$ws1.Range("M2").Formula = '=((LEFT(H2,(FIND("°",H2,1)-1)))+(MID(H2,(FIND("°",H2,1)+1),(FIND("''",H2,1))-(FIND("°",H2,1)+1))/60))*(IF(RIGHT(H2,1)="E",1,-1))'

# Style for L2/M2: centered, black Calibri 11 -- reuse the existing matching cell
# format (Asset_Cal_Info!B2) via copy/paste-special so no redundant font/style
# table entries are introduced.
$ws2.Range("B2").Copy()
$ws1.Range("L2:M2").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# --- Sheet2 (Asset_Cal_Info) corrections ---
# Fix MOPAK reference designator text
$ws2.Range("D6").Value = "CP02PMCO-MOPAK"

# Update selection on sheet2 (Asset_Cal_Info), then return focus to sheet1 (Moorings)
# so the workbook's active/selected tab ends up matching the original (Moorings).
$ws2.Range("F4").Select()
$ws1.Activate()
$ws1.Range("P9").Select()
